$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled column header: LBNDIND -> LBNRIND
$ws.Range("F1").Value = "LBNRIND"

# Update selection to reflect new active cell/range
$ws.Range("D8:D9").Select()
